# "12 - More Methods" / updates to 08-14 as part of recording session
#
# The three shapes being repositioned (the CSU footer logo, the
# "Department of Computer Science" textbox, and the credit-line
# rectangle) live on the slide layout ("Title Green Ram CSU"), not on
# the slide itself, so we reach them through the slide's CustomLayout.
#
# PowerPoint's COM Shape.Left/.Top/.Width/.Height are expressed in
# points (1 pt = 12700 EMU) and are stored in single precision, so the
# literal EMU/12700 value can truncate down by a hair once it round-trips
# through the points<->EMU conversion. The point values below are nudged
# by a hair (well under a hundredth of a point) so the stored EMU lands
# exactly on the target value.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$cl = $s.CustomLayout

# Picture 6 (footer logo image)
$picture = $cl.Shapes.Item("Picture 6")
$picture.Left   = 828.0
$picture.Top    = 527.1008661417322
$picture.Width  = 249.4652023503937
$picture.Height = 55.798425696850394

# TextBox 8 ("Department of Computer Science")
$textBox = $cl.Shapes.Item("TextBox 8")
$textBox.Left = 880.782283464567
$textBox.Top  = 563.0076599952756

# Rectangle 10 (credit-line text box)
$rectangle = $cl.Shapes.Item("Rectangle 10")
$rectangle.Left = 825.782283464567
$rectangle.Top  = 591.9631653062992
